# comment from python script
#
# Fills in the new "S3" (column D) attitude scores for every student on the
# four sheets (Astronauta, Senador, Mago, Ninja) and restores each sheet's
# selection the way it was left after the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Astronauta (sheet1): every student gets a perfect S3 score of 1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Astronauta")
$sheet1_D = @{
  2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=1; 9=1; 10=1; 11=1;
  12=1; 13=1; 14=1; 15=1; 16=1; 17=1; 18=1; 19=1; 20=1; 21=1; 22=1
}
foreach ($r in $sheet1_D.Keys) {
  $ws1.Cells.Item($r, 4).Value = $sheet1_D[$r]
}
[void]$ws1.Activate()
[void]$ws1.Range("D23").Select()

# ---------------------------------------------------------------------------
# Senador (sheet2): mixed S3 scores
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Senador")
$sheet2_D = @{
  2=1; 3=1; 4=0.4; 5=1; 6=0; 7=1; 8=1; 9=1; 10=1; 11=1;
  12=1; 13=0.5; 14=0; 15=0.5; 16=0; 17=1; 18=1; 19=1; 20=1; 21=1; 22=1
}
foreach ($r in $sheet2_D.Keys) {
  $ws2.Cells.Item($r, 4).Value = $sheet2_D[$r]
}
[void]$ws2.Activate()
[void]$ws2.Range("D5").Select()

# ---------------------------------------------------------------------------
# Mago (sheet3): S3 column already populated, only the selection moved
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Mago")
[void]$ws3.Activate()
[void]$ws3.Range("D5").Select()

# ---------------------------------------------------------------------------
# Ninja (sheet4): mixed S3 scores
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Ninja")
$sheet4_D = @{
  2=1; 3=0; 4=0; 5=1; 6=1; 7=1; 8=1; 9=1; 10=1; 11=1;
  12=1; 13=1; 14=1; 15=1; 16=1; 17=1; 18=1; 19=1; 20=1; 21=1; 22=1
}
foreach ($r in $sheet4_D.Keys) {
  $ws4.Cells.Item($r, 4).Value = $sheet4_D[$r]
}
[void]$ws4.Activate()
[void]$ws4.Range("D5").Select()

# Restore the originally active sheet (Astronauta, tab 1) and its selection
[void]$ws1.Activate()
[void]$ws1.Range("D23").Select()
